# Updates odds/score values on Sheet1 to reflect the refreshed FlashScore export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 3.25  # G2: 3.2 -> 3.25
$ws.Cells.Item(2, 9).Value = 2.25  # I2: 2.3 -> 2.25
$ws.Cells.Item(2, 17).Value = 2.1  # Q2: 2.2 -> 2.1
$ws.Cells.Item(2, 18).Value = 1.73  # R2: 1.67 -> 1.73
$ws.Cells.Item(2, 27).Value = 29  # AA2: 26 -> 29
$ws.Cells.Item(2, 28).Value = 41  # AB2: 34 -> 41
$ws.Cells.Item(2, 34).Value = 7  # AH2: 7.5 -> 7
$ws.Cells.Item(2, 35).Value = 10  # AI2: 11 -> 10
$ws.Cells.Item(2, 41).Value = 19  # AO2: 17 -> 19
$ws.Cells.Item(2, 56).Value = 151  # BD2: 126 -> 151
# Row 3
$ws.Cells.Item(3, 7).Value = 2.25  # G3: 2.35 -> 2.25
$ws.Cells.Item(3, 9).Value = 3.5  # I3: 3.3 -> 3.5
$ws.Cells.Item(3, 10).Value = 3  # J3: 3.1 -> 3
$ws.Cells.Item(3, 12).Value = 4  # L3: 3.75 -> 4
$ws.Cells.Item(3, 21).Value = 1.95  # U3: 1.91 -> 1.95
$ws.Cells.Item(3, 22).Value = 1.8  # V3: 1.91 -> 1.8
$ws.Cells.Item(3, 24).Value = 10  # X3: 11 -> 10
$ws.Cells.Item(3, 27).Value = 19  # AA3: 21 -> 19
$ws.Cells.Item(3, 29).Value = 7.5  # AC3: 8 -> 7.5
$ws.Cells.Item(3, 33).Value = 351  # AG3: 301 -> 351
$ws.Cells.Item(3, 35).Value = 17  # AI3: 15 -> 17
$ws.Cells.Item(3, 36).Value = 13  # AJ3: 12 -> 13
$ws.Cells.Item(3, 37).Value = 41  # AK3: 34 -> 41
$ws.Cells.Item(3, 40).Value = 4  # AN3: 4.33 -> 4
$ws.Cells.Item(3, 48).Value = 67  # AV3: 51 -> 67
$ws.Cells.Item(3, 53).Value = 101  # BA3: 81 -> 101
# Row 5
$ws.Cells.Item(5, 7).Value = 2.62  # G5: 2.57 -> 2.62
$ws.Cells.Item(5, 9).Value = 2.55  # I5: 2.6 -> 2.55
$ws.Cells.Item(5, 10).Value = 3.15  # J5: 3.2 -> 3.15
$ws.Cells.Item(5, 11).Value = 2.12  # K5: 2.07 -> 2.12
$ws.Cells.Item(5, 12).Value = 3.15  # L5: 3.2 -> 3.15
$ws.Cells.Item(5, 16).Value = 3.3  # P5: 3.35 -> 3.3
$ws.Cells.Item(5, 19).Value = 1.39  # S5: 1.42 -> 1.39
$ws.Cells.Item(5, 20).Value = 2.8  # T5: 2.67 -> 2.8
$ws.Cells.Item(5, 23).Value = 9.5  # W5: 9.25 -> 9.5
$ws.Cells.Item(5, 24).Value = 14  # X5: 13.5 -> 14
$ws.Cells.Item(5, 26).Value = 30  # Z5: 29 -> 30
$ws.Cells.Item(5, 34).Value = 8.75  # AH5: 9 -> 8.75
$ws.Cells.Item(5, 35).Value = 13  # AI5: 13.5 -> 13
$ws.Cells.Item(5, 37).Value = 28  # AK5: 29 -> 28
$ws.Cells.Item(5, 40).Value = 4.65  # AN5: 4.55 -> 4.65
$ws.Cells.Item(5, 41).Value = 13.5  # AO5: 14 -> 13.5
$ws.Cells.Item(5, 42).Value = 19.5  # AP5: 21 -> 19.5
$ws.Cells.Item(5, 44).Value = 80  # AR5: 90 -> 80
$ws.Cells.Item(5, 46).Value = 2.8  # AT5: 2.67 -> 2.8
$ws.Cells.Item(5, 47).Value = 6.7  # AU5: 6.9 -> 6.7
$ws.Cells.Item(5, 48).Value = 55  # AV5: 60 -> 55
$ws.Cells.Item(5, 51).Value = 21  # AY5: 22 -> 21
$ws.Cells.Item(5, 52).Value = 60  # AZ5: 65 -> 60
$ws.Cells.Item(5, 53).Value = 90  # BA5: 100 -> 90
$ws.Cells.Item(5, 54).Value = 250  # BB5: 300 -> 250
# Row 6
$ws.Cells.Item(6, 7).Value = 1.7  # G6: 1.66 -> 1.7
# Row 7
$ws.Cells.Item(7, 9).Value = 2.2  # I7: 2.15 -> 2.2
# Row 8
$ws.Cells.Item(8, 7).Value = 1.91  # G8: 1.86 -> 1.91
$ws.Cells.Item(8, 19).Value = 1.3  # S8: 1.33 -> 1.3
$ws.Cells.Item(8, 56).Value = 126  # BD8: 151 -> 126
# Row 9
$ws.Cells.Item(9, 19).Value = 1.47  # S9: 1.5 -> 1.47
# Row 11
$ws.Cells.Item(11, 13).Value = 1.06  # M11: 1.07 -> 1.06
$ws.Cells.Item(11, 14).Value = 10  # N11: 9 -> 10
# Row 16
$ws.Cells.Item(16, 7).Value = 3.75  # G16: 3.6 -> 3.75
$ws.Cells.Item(16, 9).Value = 1.95  # I16: 2 -> 1.95
$ws.Cells.Item(16, 10).Value = 4.33  # J16: 4 -> 4.33
$ws.Cells.Item(16, 11).Value = 2.2  # K16: 2.1 -> 2.2
$ws.Cells.Item(16, 12).Value = 2.63  # L16: 2.75 -> 2.63
$ws.Cells.Item(16, 13).Value = 1.06  # M16: 1.05 -> 1.06
$ws.Cells.Item(16, 14).Value = 10  # N16: 11 -> 10
$ws.Cells.Item(16, 29).Value = 10  # AC16: 9.5 -> 10
$ws.Cells.Item(16, 39).Value = 26  # AM16: 29 -> 26
$ws.Cells.Item(16, 52).Value = 34  # AZ16: 41 -> 34
# Row 20
$ws.Cells.Item(20, 7).Value = 3.9  # G20: 4 -> 3.9
$ws.Cells.Item(20, 8).Value = 3.2  # H20: 3.1 -> 3.2
$ws.Cells.Item(20, 14).Value = 9.5  # N20: 10 -> 9.5
$ws.Cells.Item(20, 27).Value = 29  # AA20: 34 -> 29
$ws.Cells.Item(20, 29).Value = 9.5  # AC20: 9 -> 9.5
$ws.Cells.Item(20, 37).Value = 19  # AK20: 17 -> 19
$ws.Cells.Item(20, 39).Value = 26  # AM20: 29 -> 26
# Row 21
$ws.Cells.Item(21, 17).Value = 2.3  # Q21: 2.35 -> 2.3
$ws.Cells.Item(21, 18).Value = 1.6  # R21: 1.57 -> 1.6
# Row 24
$ws.Cells.Item(24, 7).Value = 1.62  # G24: 1.7 -> 1.62
$ws.Cells.Item(24, 8).Value = 3.95  # H24: 3.85 -> 3.95
$ws.Cells.Item(24, 9).Value = 5  # I24: 4.6 -> 5
$ws.Cells.Item(24, 10).Value = 2.1  # J24: 2.15 -> 2.1
$ws.Cells.Item(24, 11).Value = 2.32  # K24: 2.35 -> 2.32
$ws.Cells.Item(24, 12).Value = 4.85  # L24: 4.65 -> 4.85
$ws.Cells.Item(24, 14).Value = 8.75  # N24: 8.5 -> 8.75
$ws.Cells.Item(24, 15).Value = 1.25  # O24: 1.26 -> 1.25
$ws.Cells.Item(24, 16).Value = 3.75  # P24: 3.65 -> 3.75
$ws.Cells.Item(24, 17).Value = 1.75  # Q24: 1.78 -> 1.75
$ws.Cells.Item(24, 18).Value = 2.05  # R24: 2 -> 2.05
$ws.Cells.Item(24, 20).Value = 3.2  # T24: 3.25 -> 3.2
$ws.Cells.Item(24, 24).Value = 8.25  # X24: 8.75 -> 8.25
$ws.Cells.Item(24, 26).Value = 13  # Z24: 14 -> 13
$ws.Cells.Item(24, 27).Value = 13.5  # AA24: 14 -> 13.5
$ws.Cells.Item(24, 29).Value = 8.75  # AC24: 8.5 -> 8.75
$ws.Cells.Item(24, 30).Value = 8  # AD24: 7.8 -> 8
$ws.Cells.Item(24, 34).Value = 14  # AH24: 13 -> 14
$ws.Cells.Item(24, 35).Value = 32  # AI24: 30 -> 32
$ws.Cells.Item(24, 36).Value = 17  # AJ24: 16 -> 17
$ws.Cells.Item(24, 37).Value = 100  # AK24: 90 -> 100
$ws.Cells.Item(24, 38).Value = 50  # AL24: 45 -> 50
$ws.Cells.Item(24, 40).Value = 3.6  # AN24: 3.7 -> 3.6
$ws.Cells.Item(24, 41).Value = 7.6  # AO24: 7.8 -> 7.6
$ws.Cells.Item(24, 42).Value = 15.5  # AP24: 15 -> 15.5
$ws.Cells.Item(24, 43).Value = 23  # AQ24: 24 -> 23
$ws.Cells.Item(24, 46).Value = 3.2  # AT24: 3.25 -> 3.2
$ws.Cells.Item(24, 47).Value = 7.2  # AU24: 7 -> 7.2
$ws.Cells.Item(24, 49).Value = 6.8  # AW24: 6.5 -> 6.8
$ws.Cells.Item(24, 50).Value = 26  # AX24: 24 -> 26
$ws.Cells.Item(24, 51).Value = 28  # AY24: 26 -> 28
$ws.Cells.Item(24, 52).Value = 150  # AZ24: 120 -> 150
# Row 26
$ws.Cells.Item(26, 9).Value = 1.86  # I26: 1.9 -> 1.86
